# Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) inside specific bullet /
# impact paragraphs, matching the author's commit.

$d = $word.ActiveDocument

# Word "wdColor" values are stored 0xBBGGRR (reverse byte order of the
# usual RGB hex triplet). 2C3E50 (R=2C,G=3E,B=50) -> 0x503E2C.
$metricColor = 0x503E2C

function Highlight-Metric {
    param(
        [int]$ParaIndex,   # 1-based Paragraphs.Item index
        [string[]]$Needles # substrings to bold+color, in left-to-right order
    )

    $para = $d.Paragraphs.Item($ParaIndex)
    $paraStart = $para.Range.Start
    $paraEnd = $para.Range.End

    $cursor = $paraStart
    foreach ($needle in $Needles) {
        $searchRange = $d.Range($cursor, $paraEnd)
        $found = $searchRange.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $searchRange.Font.Bold = $true
            $searchRange.Font.Color = $metricColor
            $cursor = $searchRange.End
        }
    }
}

# • Discovered systematic race coding errors ... accuracy from 23% to 64%
Highlight-Metric 9 @("23%", "64%")

# • Achieved 87% prediction accuracy ... of 71%, reducing polling error
#   margins from ±4.2% to ±2.1%
Highlight-Metric 11 @("87%", "71%", "±4.2%", "±2.1%")

# • Wrote RFP and analyzed bids from 1,200 vendors ...
Highlight-Metric 31 @("1,200")

# • Created comprehensive meta-analysis framework ... $400M ... $1B+
Highlight-Metric 46 @('$400M', '$1B')

# • Algorithm reduced mapping costs by 73.5%, saving ... $4.7M
Highlight-Metric 63 @('73.5%', '$4.7M')

# • Achieved 87% prediction accuracy ... of 71%  (KEY ACHIEVEMENTS section)
Highlight-Metric 65 @("87%", "71%")
